$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.544.00"
$ws.Range("E2").Value = "  +3.24%  "

$ws.Range("D3").Value = "1.827.58"
$ws.Range("E3").Value = "  +4.74%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.42"
$ws.Range("E5").Value = "  +2.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3829"
$ws.Range("E7").Value = "  +2.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3542"
$ws.Range("E8").Value = "  +4.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.91"
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.238"
$ws.Range("E10").Value = "  +3.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07755"
$ws.Range("E11").Value = "  +3.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.21"
$ws.Range("E13").Value = "  +8.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.613"
$ws.Range("E14").Value = "  +2.05%  "

$ws.Range("D15").Value = "1.827.56"
$ws.Range("E15").Value = "  +4.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.230"
$ws.Range("E16").Value = "  +1.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("E17").Value = "  +3.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06750"
$ws.Range("E18").Value = "  +0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.93"
$ws.Range("E19").Value = "  +3.88%  "

$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.60"
$ws.Range("E21").Value = "  +4.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.554"
$ws.Range("E22").Value = "  +5.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.19"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("D24").Value = "27.541.13"
$ws.Range("E24").Value = "  +3.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.471"
$ws.Range("E25").Value = "  +1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.712"
$ws.Range("E26").Value = "  +9.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.05"
$ws.Range("E27").Value = "  +11.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.493"
$ws.Range("E28").Value = "  +4.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.65"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D30").Value = "2.037.45"
$ws.Range("E30").Value = "  +5.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.32"
$ws.Range("E31").Value = "  +2.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.383"
$ws.Range("E32").Value = "  +4.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.091"
$ws.Range("E33").Value = "  -0.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.90"
$ws.Range("E34").Value = "  +6.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08815"
$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.641"
$ws.Range("E37").Value = "  +3.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7056"
$ws.Range("E38").Value = "  +12.54%  "

$ws.Range("E39").Value = "  +5.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06523"
$ws.Range("E40").Value = "  +3.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2261"
$ws.Range("E41").Value = "  +3.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02408"
$ws.Range("E42").Value = "  +1.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.318"
$ws.Range("E43").Value = "  +7.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.71"
$ws.Range("E44").Value = "  +3.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6623"
$ws.Range("E45").Value = "  +9.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.952"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.195"
$ws.Range("E48").Value = "  +6.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "133.52"
$ws.Range("E49").Value = "  +3.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07311"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.08"
$ws.Range("E51").Value = "  +3.83%  "
